# Weekly update: insert a new price record at row 22 (pushing the
# existing rows 22-51 down to 23-52) and populate the new row with
# this week's data for "Agrícola del Norte S.A. de Arica" / Mandarina.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 22..51 down to 23..52, leaving row 22 empty (and carrying
# the existing "D" column date-number-format down with the cells, as
# Excel's row Insert normally does).
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly record.
$ws.Cells.Item(22, 1).Value  = 1
$ws.Cells.Item(22, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(22, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(22, 4).Value  = 44413
$ws.Cells.Item(22, 5).Value  = 15
$ws.Cells.Item(22, 6).Value  = "Fruta"
$ws.Cells.Item(22, 7).Value  = 100102
$ws.Cells.Item(22, 8).Value  = "Cítricos"
$ws.Cells.Item(22, 9).Value  = 100102004
$ws.Cells.Item(22, 10).Value = "Mandarina"
$ws.Cells.Item(22, 11).Value = "Clemenuless"
$ws.Cells.Item(22, 12).Value = "Segunda"
$ws.Cells.Item(22, 13).Value = 300
$ws.Cells.Item(22, 14).Value = 12000
$ws.Cells.Item(22, 15).Value = 13000
$ws.Cells.Item(22, 16).Value = 12500
$ws.Cells.Item(22, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(22, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(22, 19).Value = 625
$ws.Cells.Item(22, 20).Value = 20
